$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.667.40"
$ws.Range("E2").Value = "  -1.06%  "
$ws.Range("D3").Value = "3.792.57"
$ws.Range("E3").Value = "  +1.34%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "3.790.58"
$ws.Range("E7").Value = "  +1.31%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("E11").Value = "  -1.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.449"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("D15").Value = "4.428.43"
$ws.Range("E15").Value = "  +1.36%  "
$ws.Range("D16").Value = "3.788.99"
$ws.Range("E16").Value = "  +1.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.47"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.54%  "
$ws.Range("D18").Value = "67.632.58"
$ws.Range("E18").Value = "  -1.10%  "
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("E21").Value = "  -4.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "457.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.95%  "
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000156"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.93%  "
$ws.Range("E27").Value = "  -2.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.45%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.30"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "29.88"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.20"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.22"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  -0.24%  "
$ws.Range("D36").Value = "3.745.62"
$ws.Range("E36").Value = "  +1.31%  "
$ws.Range("E37").Value = "  -0.92%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.34"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.25%  "
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "45.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.300"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "47.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "148.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.70%  "
$ws.Range("E49").Value = "  -4.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "389.91"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("D51").Value = "2.761.57"
$ws.Range("E51").Value = "  +2.45%  "
